$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1273
$ws1.Range("F3").Value = 665
$ws1.Range("F5").Value = 5134
$ws1.Range("F6").Value = 548
$ws1.Range("F7").Value = 9968
$ws1.Range("F8").Value = 257
$ws1.Range("F9").Value = 550
$ws1.Range("F10").Value = 99
$ws1.Range("F11").Value = 49
$ws1.Range("F12").Value = 751
$ws1.Range("F13").Value = 82

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 1

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1273
$ws4.Range("F3").Value = 665
$ws4.Range("F7").Value = 5134
$ws4.Range("F8").Value = 548
$ws4.Range("F10").Value = 9968
$ws4.Range("F11").Value = 257
$ws4.Range("F12").Value = 550
$ws4.Range("F13").Value = 99
$ws4.Range("F16").Value = 49
$ws4.Range("F17").Value = 751
$ws4.Range("F19").Value = 82
